$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44329
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 15000
$ws.Range("P2").Value = 1000

$ws.Range("D3").Value = 44313
$ws.Range("J3").Value = 40
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 14000
$ws.Range("P3").Value = 933

$ws.Range("D4").Value = 44406
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 22000
$ws.Range("L4").Value = 22000
$ws.Range("M4").Value = 22000
$ws.Range("P4").Value = 1467

$ws.Range("D5").Value = 44455
$ws.Range("J5").Value = 35
$ws.Range("K5").Value = 22000
$ws.Range("L5").Value = 22000
$ws.Range("M5").Value = 22000
$ws.Range("P5").Value = 1467

$ws.Range("D6").Value = 44309
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("P6").Value = 1000

$ws.Range("D7").Value = 44323
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 15000
$ws.Range("P7").Value = 1000

$ws.Range("D8").Value = 44341
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("P8").Value = 1000

$ws.Range("D9").Value = 44340
$ws.Range("J9").Value = 47
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 14000
$ws.Range("P9").Value = 933

$ws.Range("D10").Value = 44399
$ws.Range("J10").Value = 38
$ws.Range("K10").Value = 22000
$ws.Range("L10").Value = 22000
$ws.Range("M10").Value = 22000
$ws.Range("P10").Value = 1467

$ws.Range("D11").Value = 44327
$ws.Range("J11").Value = 35
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 15000
$ws.Range("P11").Value = 1000

$ws.Range("D12").Value = 44370
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 18000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 18000
$ws.Range("P12").Value = 1200

$ws.Range("D13").Value = 44336
$ws.Range("J13").Value = 65
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14462
$ws.Range("P13").Value = 964

$ws.Range("D14").Value = 44308
$ws.Range("J14").Value = 40
$ws.Range("K14").Value = 16000
$ws.Range("L14").Value = 16000
$ws.Range("M14").Value = 16000
$ws.Range("P14").Value = 1067

$ws.Range("D15").Value = 44330
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 15000
$ws.Range("P15").Value = 1000

$ws.Range("D16").Value = 44322
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("P16").Value = 967

$ws.Range("D17").Value = 44315
$ws.Range("J17").Value = 65
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14538
$ws.Range("P17").Value = 969

$ws.Range("D18").Value = 44343
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 15000
$ws.Range("P18").Value = 1000

$ws.Range("D19").Value = 44316
$ws.Range("J19").Value = 45
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 14444
$ws.Range("P19").Value = 963

$ws.Range("D20").Value = 44319
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 15000
$ws.Range("P20").Value = 1000

$ws.Range("D21").Value = 44411
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 22000
$ws.Range("L21").Value = 22000
$ws.Range("M21").Value = 22000
$ws.Range("P21").Value = 1467

$ws.Range("D22").Value = 44344
$ws.Range("J22").Value = 40
$ws.Range("K22").Value = 20000
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = 20000
$ws.Range("P22").Value = 1333

$ws.Range("D23").Value = 44320
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 15000
$ws.Range("P23").Value = 1000

$ws.Range("D25").Value = 44452
$ws.Range("J25").Value = 73
$ws.Range("K25").Value = 22000
$ws.Range("L25").Value = 23000
$ws.Range("M25").Value = 22479
$ws.Range("P25").Value = 1499

$ws.Range("D26").Value = 44377
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 18000
$ws.Range("L26").Value = 19000
$ws.Range("M26").Value = 18500
$ws.Range("P26").Value = 1233

$ws.Range("D27").Value = 44334
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 14000
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = 14000
$ws.Range("P27").Value = 933

$ws.Range("D28").Value = 44328
$ws.Range("J28").Value = 38
$ws.Range("K28").Value = 15000
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 15000
$ws.Range("P28").Value = 1000

$ws.Range("D29").Value = 44321
$ws.Range("J29").Value = 38
$ws.Range("K29").Value = 15000
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = 15000
$ws.Range("P29").Value = 1000

$ws.Range("D30").Value = 44397
$ws.Range("J30").Value = 73
$ws.Range("K30").Value = 21000
$ws.Range("L30").Value = 22000
$ws.Range("M30").Value = 21521
$ws.Range("P30").Value = 1435

$ws.Range("D31").Value = 44333
$ws.Range("J31").Value = 35
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 15000
$ws.Range("P31").Value = 1000

$ws.Range("D32").Value = 44314
$ws.Range("J32").Value = 45
$ws.Range("K32").Value = 15000
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = 15000
$ws.Range("P32").Value = 1000

$ws.Range("D34").Value = 44312
$ws.Range("J34").Value = 80
$ws.Range("K34").Value = 13000
$ws.Range("L34").Value = 14000
$ws.Range("M34").Value = 13562
$ws.Range("P34").Value = 904

$ws.Range("D35").Value = 44326
$ws.Range("J35").Value = 45
$ws.Range("K35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = 15000
$ws.Range("P35").Value = 1000
